$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh re-sorted the price rows: each destination row (2-32)
# now holds the data that used to live in a different row. Row 31 keeps its
# own data (maps to itself). Columns A-T are copied in full per row so the
# whole record (date, quality, volume, prices, unit, origin, etc.) moves
# together.
#
# destination row -> source row (as it was before this edit)
$rowMap = @{}
$rowMap[2] = 6
$rowMap[3] = 18
$rowMap[4] = 11
$rowMap[5] = 12
$rowMap[6] = 15
$rowMap[7] = 20
$rowMap[8] = 22
$rowMap[9] = 27
$rowMap[10] = 28
$rowMap[11] = 14
$rowMap[12] = 16
$rowMap[13] = 9
$rowMap[14] = 10
$rowMap[15] = 2
$rowMap[16] = 21
$rowMap[17] = 30
$rowMap[18] = 5
$rowMap[19] = 3
$rowMap[20] = 17
$rowMap[21] = 4
$rowMap[22] = 13
$rowMap[23] = 7
$rowMap[24] = 8
$rowMap[25] = 32
$rowMap[26] = 29
$rowMap[27] = 23
$rowMap[28] = 24
$rowMap[29] = 25
$rowMap[30] = 26
$rowMap[31] = 31
$rowMap[32] = 19

$firstCol = 1   # A
$lastCol  = 20  # T

# 1) Snapshot every source row BEFORE any writes (the mapping is a full
#    permutation, so a naive in-place copy would clobber data that is still
#    needed as a source for another destination row).
$snapshot = @{}
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowValues = @{}
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value()
        }
        $snapshot[$srcRow] = $rowValues
    }
}

# 2) Write the snapshotted source row into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowValues = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $rowValues[$col]
    }
}
